$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.294434905052185
$ws.Range("B1").Value = 1.994203567504883
$ws.Range("C1").Value = 5.326944351196289
$ws.Range("D1").Value = 1.928382396697998
$ws.Range("E1").Value = 1.095526576042175
